$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new log entries (rows 11-13)
$ws.Range("A11").Value = "alleen"
$ws.Range("B11").Value = "bugs gateway oplossen"
$ws.Range("C11").Value = 43402
$ws.Range("D11").Value = 0.541666666666667
$ws.Range("E11").Value = 0.625

$ws.Range("A12").Value = "alleen"
$ws.Range("B12").Value = "bugs gateway oplossen"
$ws.Range("C12").Value = 43403
$ws.Range("D12").Value = 0.5625
$ws.Range("E12").Value = 0.666666666666667

$rsquo = [char]0x2019
$ws.Range("A13").Value = "alleen"
$ws.Range("B13").Value = "meerdere arduino" + $rsquo + "s op de gateway"
$ws.Range("C13").Value = 43404
$ws.Range("D13").Value = 0.5625
$ws.Range("E13").Value = 0.6875

$ws.Range("C11:C13").NumberFormat = "DD/MM/YY"
$ws.Range("D11:E13").NumberFormat = "HH:MM:SS"

$ws.Range("E34").Select()
